# Update provincias Spain data on the "Ciudades" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Mayo de 2020 a las 12:05"

# Row 4 - Madrid
$ws.Range("B4").Value = 66338
$ws.Range("C4").Value = 40698
$ws.Range("D4").Value = 16793
$ws.Range("E4").Value = 8847

# Row 5 - Cataluña
$ws.Range("B5").Value = 55824
$ws.Range("C5").Value = 26172
$ws.Range("D5").Value = 23708
$ws.Range("E5").Value = 5944

# Row 6 - Castilla y Leon
$ws.Range("B6").Value = 18470
$ws.Range("C6").Value = 8716
$ws.Range("D6").Value = 7808
$ws.Range("E6").Value = 1946

# Row 7 - Castilla-La Mancha
$ws.Range("B7").Value = 16618
$ws.Range("C7").Value = 6378
$ws.Range("D7").Value = 7347
$ws.Range("E7").Value = 2893

# Row 9 - Andalucia
$ws.Range("B9").Value = 12450
$ws.Range("C9").Value = 10611
$ws.Range("D9").Value = 481
$ws.Range("E9").Value = 1358

# Row 14 - Aragon
$ws.Range("B14").Value = 5478
$ws.Range("C14").Value = 3727
$ws.Range("D14").Value = 913

# Row 16 - Navarra
$ws.Range("B16").Value = 5148
$ws.Range("C16").Value = 3751
$ws.Range("D16").Value = 894
$ws.Range("E16").Value = 503

# Row 20 - La Rioja
$ws.Range("B20").Value = 4024
$ws.Range("C20").Value = 3048
$ws.Range("D20").Value = 627
$ws.Range("E20").Value = 349

# Row 32 - Asturias
$ws.Range("C32").Value = 1061
$ws.Range("D32").Value = 990
$ws.Range("E32").Value = 315

# Row 33 - Gran Canaria
$ws.Range("B33").Value = 2289
$ws.Range("D33").Value = 614
